# Update the score sheet with results from "ejecucion_2" (base0_ci_rf_st.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "_tejgkft_redr"
$ws.Range("C2").Value = 0.01490817827358855
$ws.Range("B3").Value = "_tejgct_r09gstcp"
$ws.Range("C3").Value = 0.0124981708024734
$ws.Range("B4").Value = "_tejgfun_f5r18ct05pgrco"
$ws.Range("C4").Value = 0.009504508568790811
$ws.Range("B5").Value = "_tejgge_r09ct06acanf"
$ws.Range("C5").Value = 0.009235829692581639
$ws.Range("B6").Value = "_devppimtotfun_f1trans"
$ws.Range("C6").Value = 0.008500909050737682
$ws.Range("B7").Value = "_tejgfun_f5ct05prots"
$ws.Range("C7").Value = 0.007437263570595411
$ws.Range("B8").Value = "_tejgfun_f5r08ct05amb"
$ws.Range("C8").Value = 0.004953404496410972
$ws.Range("B9").Value = "devppimfun_f5r07ct05salud"
$ws.Range("C9").Value = 0.004467190497387235
$ws.Range("B10").Value = "_tejgtotfun_f5r08amb"
$ws.Range("C10").Value = 0.004200475626676016
$ws.Range("B11").Value = "_tejgfun_f5ct06amb"
$ws.Range("C11").Value = 0.004021918504413352
$ws.Range("B12").Value = "_tejgfun_f5ct06opseg"
$ws.Range("C12").Value = 0.004001136126687925
$ws.Range("B13").Value = "_tejgfun_f5ct06cydep"
$ws.Range("C13").Value = 0.003651256129579302
$ws.Range("B14").Value = "tejgtotfun_f2prots"
$ws.Range("C14").Value = 0.003120078556004896
$ws.Range("B15").Value = "tejgtotfun_f5r08ambpc"
$ws.Range("C15").Value = 0.002914618686877528
$ws.Range("B16").Value = "_tejgtotfun_f2opsegpc"
$ws.Range("C16").Value = 0.00285483929644664
$ws.Range("B17").Value = "_tejgfun_f5r08ct05prots"
$ws.Range("C17").Value = 0.002760516310882607
$ws.Range("B18").Value = "tejgfun_f2ct05prots"
$ws.Range("C18").Value = 0.002596578198346151
$ws.Range("B19").Value = "_tejgtotfun_f5r08pgrco"
$ws.Range("C19").Value = 0.002514073355859688
$ws.Range("B20").Value = "_tejgtotfun_f5cydep"
$ws.Range("C20").Value = 0.002490966968799659
$ws.Range("B21").Value = "pimgtotfun_f1san"
$ws.Range("C21").Value = 0.002361647571709484
$ws.Range("B22").Value = "_tejgtotfun_f5r18amb"
$ws.Range("C22").Value = 0.002343364185036436
$ws.Range("B23").Value = "devppimtotfun_f5r07salud"
$ws.Range("C23").Value = 0.002275141470269209
$ws.Range("B24").Value = "tejgtotfun_f5r18opseg"
$ws.Range("C24").Value = 0.002254442895210739
$ws.Range("B25").Value = "tejgfun_f5r08ct05ambpc"
$ws.Range("C25").Value = 0.002221589653770451
$ws.Range("B26").Value = "devppimfun_f5r07ct05agro"
$ws.Range("C26").Value = 0.002196389907218644
$ws.Range("B27").Value = "pimgfun_f5r18ct05trans"
$ws.Range("C27").Value = 0.002122591581612303
$ws.Range("B28").Value = "tejgfun_f5ct06viv"
$ws.Range("C28").Value = 0.00203916844221394
$ws.Range("B29").Value = "_pimgfun_f5ct06opsegpc"
$ws.Range("C29").Value = 0.001986023024234364
$ws.Range("B30").Value = "_devppimfun_f1ct05trans"
$ws.Range("C30").Value = 0.001973829720159832
$ws.Range("B31").Value = "_devppimtotfun_f5viv"
$ws.Range("C31").Value = 0.001963577945787877
$ws.Range("B32").Value = "_tejgtotfun_f5r18cydep"
$ws.Range("C32").Value = 0.001942393940575374
$ws.Range("B33").Value = "_tejgfun_f5ct06opsegpc"
$ws.Range("C33").Value = 0.001915156092049192
$ws.Range("B34").Value = "tejgfun_f5ct05prots"
$ws.Range("C34").Value = 0.001823239698289403
$ws.Range("B35").Value = "pimgge_r00ct06acanfpc"
$ws.Range("C35").Value = 0.001821775709418529
$ws.Range("B36").Value = "_pimgfun_f5r18ct06opseg"
$ws.Range("C36").Value = 0.001814213266710174
$ws.Range("B37").Value = "pimgfun_f1ct06san"
$ws.Range("C37").Value = 0.001780730405485892
$ws.Range("B38").Value = "dfgpimpiafun_f1ct05prots"
$ws.Range("C38").Value = 0.001775965807087391
$ws.Range("B39").Value = "_tejgge_r08ct05biser"
$ws.Range("C39").Value = 0.001741090494929363
$ws.Range("B40").Value = "_dfgdevpiagrb_foncpc"
$ws.Range("C40").Value = 0.001621785359653043
$ws.Range("B41").Value = "_tejgtotfun_f5opseg"
$ws.Range("C41").Value = 0.001611630791843809
$ws.Range("B42").Value = "pimgkft_reodpc"
$ws.Range("C42").Value = 0.001502775324889154
$ws.Range("B43").Value = "pimgct_r00gstcppc"
$ws.Range("C43").Value = 0.001487750372268488
$ws.Range("B44").Value = "_tejgct_r00gstcrpc"
$ws.Range("C44").Value = 0.001484359269316874
$ws.Range("B45").Value = "_tejgtotfun_f2opseg"
$ws.Range("C45").Value = 0.001443890192960734
$ws.Range("B46").Value = "devppimtotfun_f5turi"
$ws.Range("C46").Value = 0.001413312387911662
$ws.Range("B47").Value = "dfgdevpiagfun_f5ct05sanpc"
$ws.Range("C47").Value = 0.001401774794793889
$ws.Range("B48").Value = "_tejgfun_f1ct05protspc"
$ws.Range("C48").Value = 0.00137420752723748
$ws.Range("B49").Value = "_devppimfun_f5r18ct06agro"
$ws.Range("C49").Value = 0.001370652468550006
$ws.Range("B50").Value = "_tdvgfun_f5ct06opsegpc"
$ws.Range("C50").Value = 0.001370254284789871
$ws.Range("B51").Value = "_tejgrb_impm"
$ws.Range("C51").Value = 0.001367869251908179
